# Refactor open account test using page objects model
#
# - The "OpenAccountTest" page-object now drives the customer-name field,
#   so the sample test data is updated from "Joao Silva" to "Harry Potter".
# - The active sheet/selection moves from "AddCustomerTest" (A1) to
#   "OpenAccountTest", with the cursor left on F9 (the new field the page
#   object interacts with).

$wb = $excel.ActiveWorkbook

$addCustomer  = $wb.Worksheets.Item("AddCustomerTest")
$openAccount  = $wb.Worksheets.Item("OpenAccountTest")

# Update the sample customer name used by the OpenAccountTest test data.
$openAccount.Range("A2").Value = "Harry Potter"

# Make sure the previously active sheet's own selection is left untouched
# before we move focus elsewhere.
$addCustomer.Range("E7").Select()

# Switch the active sheet/tab to OpenAccountTest and park the selection on
# F9, matching the page object's new entry point.
$openAccount.Activate()
$openAccount.Range("F9").Select()
